$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Hide rows 45-47 (they become hidden in the revised log)
# ---------------------------------------------------------------------
$ws.Rows("45:47").Hidden = $true

# ---------------------------------------------------------------------
# 2. Resize column A and column B
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.428571428571429
$ws.Columns.Item(2).ColumnWidth = 14.714285714285715

# ---------------------------------------------------------------------
# 3. Zoom the view to 70%
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70

# ---------------------------------------------------------------------
# 4. Add the new work-log entry as row 51, before touching row 50's
#    own formatting, so we can reuse row 50's current (pre-edit) look
#    for the freshly appended row.
# ---------------------------------------------------------------------
$ws.Range("A51").Value = 45757
$ws.Range("B51").Value = 2
$ws.Range("C51").Value = "Rehearsed with the laptop borrowed from school library, continue generating final report"

$ws.Range("A50:C50").Copy() | Out-Null
$ws.Range("A51:C51").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 5. Row 50 reverts to the regular (non-"last row") look, matching the
#    rest of the table -- copy formatting down from row 49.
# ---------------------------------------------------------------------
$ws.Range("A49:C49").Copy() | Out-Null
$ws.Range("A50:C50").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 6. Restore the selection to the new bottom of the log
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("C56").Select() | Out-Null
